# Insert a new data row above row 166 (pushing the existing row 166..259
# block down to 167..260) and populate the newly-inserted row 166 with a
# fresh observation. All columns except D (Fecha), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado) and
# P (Precio $/Kg) are identical to the row that is now at 167, so the
# simplest and most faithful way to reproduce this is:
#   1. Insert a blank row at position 166 (native row-insert semantics:
#      formatting/row height of row 166 carries down with the shifted data).
#   2. Copy the (now shifted) row 167 values into the new row 166.
#   3. Overwrite the six changed cells with their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(166).Insert()

$templateRow = $ws.Range("A167:R167").Value()
$ws.Range("A166:R166").Value = $templateRow

$ws.Range("D166").Value = 44572
$ws.Range("J166").Value = 170
$ws.Range("K166").Value = 4000
$ws.Range("L166").Value = 4500
$ws.Range("M166").Value = 4235
$ws.Range("P166").Value = 1412
